$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $rng = $d.Content
    $any = $false
    $found = $rng.Find.Execute($find, $true, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
    while ($found) {
        $rng.Text = $replace
        $any = $true
        $rng.Collapse(0)
        $rng.SetRange($rng.End, $d.Content.End)
        $found = $rng.Find.Execute($find, $true, $false, $false, $false, $false, `
                                    $true, 1, $false, "", 0)
    }
    return $any
}

Replace-Text "The playful mathematicians:" "Wanahisabati wanaocheza:"
Replace-Text "** the dialogue starts at second 47, so I added 28 seconds to all the times as they were. -John Argentino" "** mazungumzo huanza saa 47 ya pili, kwa hivyo niliongeza sekunde 28 kwa nyakati zote kama zilivyokuwa. -John Argentino"
Replace-Text "[Music]" "[Muziki]"
Replace-Text "there are two mathematicians, let's call" "kuna wanahisabati wawili, tupige simu"
Replace-Text "them Fil and Mike who meet each other" "Fil na Mike wanaokutana"
Replace-Text "again after a long time. After some" "tena baada ya muda mrefu. Baada ya baadhi"
Replace-Text "chatting, Phil says he has three children, then" "kuzungumza, Phil anasema ana watoto watatu, basi"
Replace-Text "Mike, astonished, asks: 'How old are they?' Fil," "Kwa mshangao, Mike anauliza: 'Wana umri gani?' Fil,"
Replace-Text "being a playful mathematician, answers" "kuwa mwanahisabati mchezaji, anajibu"
Replace-Text "'You tell me! I'll give you a hint: if you" "'Wewe niambie! Nitakupa kidokezo: ikiwa wewe"
Replace-Text "multiply the three ages together you" "zidisheni enzi tatu pamoja ninyi"
Replace-Text "get 36.' Mike takes sometimes to think" "pata 36.' Mike huchukua wakati mwingine kufikiria"
Replace-Text "and says: 'I'm sorry Fil, but I do need" "na kusema: 'Samahani Fil, lakini nahitaji"
Replace-Text "another hint. So Fil tells Mike:" "kidokezo kingine. Kwa hivyo Fil anamwambia Mike:"
Replace-Text "'Yes, sure, here it is: if you had up to" "'Ndiyo, hakika, hapa ni: kama alikuwa na hadi"
Replace-Text "three ages you get the number of math" "miaka mitatu unapata idadi ya hesabu"
Replace-Text "papers we publish together. Do you remember it?'" "karatasi tunachapisha pamoja. Je, unaikumbuka?'"
Replace-Text "'Yes I do remember How many, but still" "'Ndio nakumbuka wangapi, lakini bado"
Replace-Text "I do not have enough information! I need" "Sina taarifa za kutosha! nahitaji"
Replace-Text "at least one more.' Fil says: 'Yes don't" "angalau moja zaidi.' Fil anasema: 'Ndiyo usifanye hivyo"
Replace-Text "worry but this is the last one:" "wasiwasi lakini hii ni ya mwisho:"
Replace-Text "The youngest one has blues eyes.' And" "Mdogo ana macho ya blues.' Na"
Replace-Text "suddenly Mike gets the answer. You" "ghafla Mike anapata jibu. Wewe"
Replace-Text "hear the conversation but you don't know" "sikia mazungumzo lakini hujui"
Replace-Text "how many papers they published together." "ni karatasi ngapi walichapisha pamoja."
Replace-Text "However, you do want to know the ages of" "Hata hivyo, unataka kujua umri wa"
Replace-Text "the three children. Can you figure them" "watoto watatu. Je, unaweza kuwahesabu"
Replace-Text "out?" "nje?"
